$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  B=995618.8601498967;  C=61354.54598056434},
    @{Row=3;  B=515260.7276278729;  C=61354.54598056434},
    @{Row=4;  B=320071.2597331816;  C=42139.2552890066},
    @{Row=5;  B=210833.5543273592;  C=42139.2552890066},
    @{Row=6;  B=150656.9757443959;  C=42139.2552890066},
    @{Row=7;  B=107715.6238173954;  C=44260.27809068429},
    @{Row=8;  B=88064.06836449282;  C=42251.41380406159},
    @{Row=9;  B=84319.26017831972;  C=48844.82152152419},
    @{Row=10; B=79648.7828545571;   C=41811.85909840682},
    @{Row=11; B=69086.08715108788;  C=41811.85909840682},
    @{Row=12; B=75837.6307985358;   C=41811.85909840682},
    @{Row=13; B=74369.97082937368;  C=39591.65548280146},
    @{Row=14; B=66414.48470290873;  C=39591.65548280146},
    @{Row=15; B=64826.18311801473;  C=39591.65548280146},
    @{Row=16; B=72735.32619288626;  C=39591.65548280146},
    @{Row=17; B=64412.8942980627;   C=37579.2200134087},
    @{Row=18; B=59864.39898390236;  C=39591.65548280146},
    @{Row=19; B=64809.4527799048;   C=39591.65548280146},
    @{Row=20; B=67812.26166860352;  C=39591.65548280146},
    @{Row=21; B=63441.46508411197;  C=36852.7428720292}
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 2).Value = $item.B
    $ws.Cells.Item($item.Row, 3).Value = $item.C
}
